$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 17 (this pushes the old row 17..19 summary rows down to 18..20)
$ws.Rows.Item(17).Insert()

# Fill in the previously-empty row 16 with a new time entry
$ws.Cells.Item(16, 1).Value = 2014
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(16, 3).Value = 22
$ws.Cells.Item(16, 4).Value = 0.35416666666666669
$ws.Cells.Item(16, 5).Value = 0.39583333333333331
$ws.Range("F16").Formula = "=(E16-D16)*24*60"
$ws.Range("G16").Formula = "=F16/60"

# Update the summary formulas (now on rows 18-20) to cover the extended data range
$ws.Range("F18").Formula = "=SUM(F2:F17)"
$ws.Range("F19").Formula = "=F18/60"
$ws.Range("F20").Formula = "=F19/38.5"

# Update the active selection
$ws.Range("F16").Select()
